# "Generate Report for handoff"
#
# b.md.md has been handed off again (new source revision
# b3a40d6229ff1a8b48804fcfc66c95922eb78fd0). Update the Overview sheet and
# the per-locale (zh-cn / de-de) detail sheets so row 3 (b.md.md) reflects
# the new "Ready for handoff" status together with the freshly generated
# handoff file names / timestamps. Row 2 (a.md.md) and row 4
# (.localization-config) are untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-14 03:06:01"

# --- de-de detail sheet -----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-14 03:06:14"
